$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - first sheet
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 2894
$wsExhibit.Range("F6").Value = 609

# Sheet "全部类型" (All types) - fourth sheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 2894
$wsAll.Range("F8").Value = 609
